$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.897.85"
$ws.Range("E2").Value = "  -2.71%  "
$ws.Range("D3").Value = "3.865.74"
$ws.Range("E3").Value = "  -2.76%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.80"
$ws.Range("E5").Value = "  +0.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.26"
$ws.Range("E6").Value = "  +4.64%  "
$ws.Range("E7").Value = "  -1.69%  "
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("E9").Value = "  +0.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.175"
$ws.Range("E10").Value = "  +4.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.46"
$ws.Range("E11").Value = "  -1.12%  "
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.48"
$ws.Range("E13").Value = "  +5.49%  "
$ws.Range("D14").Value = "4.493.47"
$ws.Range("E14").Value = "  -2.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.38"
$ws.Range("E15").Value = "  +5.16%  "
$ws.Range("D16").Value = "3.885.62"
$ws.Range("E16").Value = "  -2.21%  "
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("E18").Value = "  -4.22%  "
$ws.Range("E19").Value = "  -2.11%  "
$ws.Range("D20").Value = "70.894.63"
$ws.Range("E20").Value = "  -2.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "436.99"
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.73"
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "94.33"
$ws.Range("E23").Value = "  -1.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.28"
$ws.Range("E24").Value = "  -4.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.92"
$ws.Range("E25").Value = "  -2.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.38"
$ws.Range("E26").Value = "  +2.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.07"
$ws.Range("E27").Value = "  -7.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.93"
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("E29").Value = "  -1.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.17"
$ws.Range("E30").Value = "  -3.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.13"
$ws.Range("E31").Value = "  +3.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.58"
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "48.09"
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("E34").Value = "  -4.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "70.03"
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("D36").Value = "0.0₃0984"
$ws.Range("E36").Value = "  +12.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "634.01"
$ws.Range("E37").Value = "  -6.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.427"
$ws.Range("E38").Value = "  -1.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.147"
$ws.Range("E39").Value = "  +0.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.31"
$ws.Range("E41").Value = "  -2.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.29"
$ws.Range("E43").Value = "  +27.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0471"
$ws.Range("E44").Value = "  -3.34%  "
$ws.Range("E45").Value = "  -7.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.72"
$ws.Range("E46").Value = "  +2.62%  "
$ws.Range("E47").Value = "  -3.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.33"
$ws.Range("E48").Value = "  -2.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.83"
$ws.Range("E49").Value = "  -15.16%  "
$ws.Range("D50").Value = "2.846.82"
$ws.Range("E50").Value = "  +1.26%  "
$ws.Range("E51").Value = "  +1.38%  "
